# The weekly data export for this sheet moved on by one week: a brand-new
# pair of rows (Primera/Segunda) for the latest date is prepended at the
# top of the data block (row 148), pushing the previously-top entry (and
# everything below it) down by two rows. The used range grows from R219 to
# R221 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (rows 148-149);
# everything currently at row 148 and below shifts down to 150+.
$ws.Range("148:149").Insert()

# New row 148 - "Primera" quality, newest date in the series.
$ws.Range("A148").Value = 8
$ws.Range("B148").Value = "Terminal La Palmera de La Serena"
$ws.Range("C148").Value = "Coquimbo"
$ws.Range("D148").Value = 44452
$ws.Range("E148").Value = 4
$ws.Range("F148").Value = 100112009
$ws.Range("G148").Value = "Acelga"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 2900
$ws.Range("K148").Value = 450
$ws.Range("L148").Value = 500
$ws.Range("M148").Value = 475
$ws.Range("N148").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O148").Value = "Provincia del Elquí"
$ws.Range("P148").Value = 238
$ws.Range("Q148").Value = 2
$ws.Range("R148").Value = "Hortaliza"

# New row 149 - "Segunda" quality, newest date in the series.
$ws.Range("A149").Value = 8
$ws.Range("B149").Value = "Terminal La Palmera de La Serena"
$ws.Range("C149").Value = "Coquimbo"
$ws.Range("D149").Value = 44452
$ws.Range("E149").Value = 4
$ws.Range("F149").Value = 100112009
$ws.Range("G149").Value = "Acelga"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Segunda"
$ws.Range("J149").Value = 1400
$ws.Range("K149").Value = 350
$ws.Range("L149").Value = 400
$ws.Range("M149").Value = 375
$ws.Range("N149").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O149").Value = "Provincia del Elquí"
$ws.Range("P149").Value = 188
$ws.Range("Q149").Value = 2
$ws.Range("R149").Value = "Hortaliza"

# Keep the date column's display format consistent with the rest of
# column D (YYYY-MM-DD HH:MM:SS, style index 2 in the original file).
$ws.Range("D148:D149").NumberFormat = "YYYY-MM-DD HH:MM:SS"
